# "Generate Report for Handback"
#
# For each localized-language sheet (zh-cn, de-de) this:
#   - marks rows 2/3 as handed back (Status column, C)
#   - stamps the handback completion datetime (column H)
#   - records the latest target file / latest handback file
#     (new columns F/G) mirroring the source file (A) and the
#     handoff file (D) links, each with their own hyperlink.
# The "Ready for handoff" status label is shared (same underlying
# string) with the Overview sheet's per-language status cells, so
# those are updated too.

function Get-HyperlinkAddress {
    param($ws, $rangeAddr)
    foreach ($hl in $ws.Hyperlinks) {
        if ($hl.Range.Address() -eq $rangeAddr) {
            return $hl.Address
        }
    }
    return $null
}

$wb = $excel.ActiveWorkbook

$statusText = "Handed back: in sync with en-US"

$sheets = @(
    @{ Name = "zh-cn"; HandbackTime = "2016-03-25 11:26:18" },
    @{ Name = "de-de"; HandbackTime = "2016-03-25 11:26:34" }
)

foreach ($info in $sheets) {
    $ws = $wb.Worksheets.Item($info.Name)

    foreach ($row in 2,3) {
        # Status column (C): "Ready for handoff" -> handed back
        $ws.Cells.Item($row, 3).Value = $statusText

        # Latest Handback DateTime column (H)
        $ws.Cells.Item($row, 8).Value = $info.HandbackTime

        # Mirror the source file hyperlink (A) into the
        # "Latest Target File" column (F)
        $srcAddr = $ws.Cells.Item($row, 1).Address()
        $srcUrl = Get-HyperlinkAddress $ws $srcAddr
        $srcDisplay = $ws.Cells.Item($row, 1).Value2
        $ws.Hyperlinks.Add($ws.Cells.Item($row, 6), $srcUrl, [System.Type]::Missing, [System.Type]::Missing, $srcDisplay)

        # Mirror the handoff-file hyperlink (D) into the
        # "Latest Handback File" column (G)
        $handoffAddr = $ws.Cells.Item($row, 4).Address()
        $handoffUrl = Get-HyperlinkAddress $ws $handoffAddr
        $handoffDisplay = $ws.Cells.Item($row, 4).Value2
        $ws.Hyperlinks.Add($ws.Cells.Item($row, 7), $handoffUrl, [System.Type]::Missing, [System.Type]::Missing, $handoffDisplay)
    }
}

# Overview sheet: columns B (zh-cn) and C (de-de) show the same
# status text per source-file row.
$overview = $wb.Worksheets.Item("Overview")
foreach ($row in 2,3) {
    $overview.Cells.Item($row, 2).Value = $statusText
    $overview.Cells.Item($row, 3).Value = $statusText
}
